$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Bmp7"
$ws.Cells.Item(2, 3).Value = "Bmpr1b"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.668521
$ws.Cells.Item(2, 8).Value = 5.005563
$ws.Cells.Item(2, 9).Value = 0.9677024783929865
$ws.Cells.Item(2, 10).Value = 0.9677024783929865
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.57938
$ws.Cells.Item(2, 14).Value = 4.73814
$ws.Cells.Item(2, 15).Value = 0.7235635290775982
$ws.Cells.Item(2, 16).Value = 0.7235635290775982
$ws.Cells.Item(2, 17).Value = 2.63522869698
$ws.Cells.Item(2, 18).Value = 23.71705827282
$ws.Cells.Item(2, 19).Value = 0.7001942203631676
$ws.Cells.Item(2, 20).Value = 0.7001942203631676

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Bmp7"
$ws.Cells.Item(3, 3).Value = "Bmpr1b"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.668521
$ws.Cells.Item(3, 8).Value = 5.005563
$ws.Cells.Item(3, 9).Value = 0.9677024783929865
$ws.Cells.Item(3, 10).Value = 0.9677024783929865
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.6034
$ws.Cells.Item(3, 14).Value = 1.8102
$ws.Cells.Item(3, 15).Value = 0.2764364709224018
$ws.Cells.Item(3, 16).Value = 0.2764364709224018
$ws.Cells.Item(3, 17).Value = 1.0067855714
$ws.Cells.Item(3, 18).Value = 9.0610701426
$ws.Cells.Item(3, 19).Value = 0.2675082580298189
$ws.Cells.Item(3, 20).Value = 0.2675082580298189

# Row 4
$ws.Cells.Item(4, 1).Value = "sCs"
$ws.Cells.Item(4, 2).Value = "Bmp7"
$ws.Cells.Item(4, 3).Value = "Bmpr1b"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.05568766666666666
$ws.Cells.Item(4, 8).Value = 0.167063
$ws.Cells.Item(4, 9).Value = 0.03229752160701353
$ws.Cells.Item(4, 10).Value = 0.03229752160701353
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 1.57938
$ws.Cells.Item(4, 14).Value = 4.73814
$ws.Cells.Item(4, 15).Value = 0.7235635290775982
$ws.Cells.Item(4, 16).Value = 0.7235635290775982
$ws.Cells.Item(4, 17).Value = 0.08795198697999998
$ws.Cells.Item(4, 18).Value = 0.7915678828199999
$ws.Cells.Item(4, 19).Value = 0.02336930871443069
$ws.Cells.Item(4, 20).Value = 0.02336930871443069

# Row 5
$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Bmp7"
$ws.Cells.Item(5, 3).Value = "Bmpr1b"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.05568766666666666
$ws.Cells.Item(5, 8).Value = 0.167063
$ws.Cells.Item(5, 9).Value = 0.03229752160701353
$ws.Cells.Item(5, 10).Value = 0.03229752160701353
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.6034
$ws.Cells.Item(5, 14).Value = 1.8102
$ws.Cells.Item(5, 15).Value = 0.2764364709224018
$ws.Cells.Item(5, 16).Value = 0.2764364709224018
$ws.Cells.Item(5, 17).Value = 0.03360193806666667
$ws.Cells.Item(5, 18).Value = 0.3024174426
$ws.Cells.Item(5, 19).Value = 0.008928212892582839
$ws.Cells.Item(5, 20).Value = 0.008928212892582839

